$wb = $excel.ActiveWorkbook

# --- Update existing data values ---

# Books sheet: D4 10 -> 7
$booksWs = $wb.Worksheets.Item("Books")
$booksWs.Range("D4").Value = 7

# Sports Kit sheet: D4 10 -> 8
$sportsWs = $wb.Worksheets.Item("Sports Kit")
$sportsWs.Range("D4").Value = 8

# --- Add new "Discount" worksheet at the end ---

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$discountWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$discountWs.Name = "Discount"

$discountWs.Range("A1").Value = "Today's Discount"
$discountWs.Range("B1").Value = 10
$discountWs.Range("C1").Value = "%"

$discountWs.Columns.Item(1).ColumnWidth = 15.3

$discountWs.Range("B1").Select()

$discountWs.Activate()
